$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Samples tab's SQL query (cell B3): the query was reworked to
# drop the "Tumor" / "Analyte Type" columns from the SELECT list, keeping
# everything else (joins, WHERE, ORDER BY, LIMIT) identical.
$newSamplesQuery = @"
SELECT
    DISTINCT (smp.sample_id) AS "Sample ID",
    sp.participant_id AS "Participant ID", 
    s.study_name AS "Study Name",
    s.phs_accession AS Accession
FROM 
    df_participant sp
JOIN 
    df_study s ON sp."study.phs_accession" = s.phs_accession
JOIN 
    df_sample smp ON smp."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_diagnosis d ON d."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_program p ON p.program_acronym = s."program.program_acronym"
JOIN
    df_file f1 ON f1."sample.sample_id" = smp.sample_id
JOIN
    df_genomic_info gi ON gi."file.file_id" = f1.file_id
WHERE 
   s.phs_accession = 'phs001819' AND smp.sample_type = 'Blood Derived Normal'
ORDER BY 
    smp.sample_id ASC
LIMIT 100;
"@

# The here-string above adds a trailing newline; the original cell text has
# none, so trim it back off before writing the value.
$newSamplesQuery = $newSamplesQuery.TrimEnd("`r", "`n")

$ws.Range("B3").Value = $newSamplesQuery

# Restore the view/selection to reflect the edited cell (row 3 instead of row 4)
$ws.Range("C3").Select()
